# Fruta / hortaliza, semanal
# New weekly price observation is inserted as a new row 213 on the
# "Hortaliza, Vega Modelo de Temuco - Ciboulette" sheet, pushing the
# existing rows 213:295 down to 214:296 (dimension grows from A1:R295 to
# A1:R296).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 213 (shifts 213:295 -> 214:296).
$ws.Rows("213:213").Insert()

# Populate the new row with the latest weekly observation.
$ws.Cells.Item(213, 1).Value  = 10
$ws.Cells.Item(213, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(213, 3).Value  = "La Araucanía"
$ws.Cells.Item(213, 4).Value  = 44875
$ws.Cells.Item(213, 5).Value  = 9
$ws.Cells.Item(213, 6).Value  = 100112039
$ws.Cells.Item(213, 7).Value  = "Ciboulette"
$ws.Cells.Item(213, 8).Value  = "Sin especificar"
$ws.Cells.Item(213, 9).Value  = "Primera"
$ws.Cells.Item(213, 10).Value = 85
$ws.Cells.Item(213, 11).Value = 6000
$ws.Cells.Item(213, 12).Value = 6000
$ws.Cells.Item(213, 13).Value = 6000
$ws.Cells.Item(213, 14).Value = "`$/docena de atados"
$ws.Cells.Item(213, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(213, 16).Value = 2000
$ws.Cells.Item(213, 17).Value = 3
$ws.Cells.Item(213, 18).Value = "Hortaliza"
